$wb = $excel.ActiveWorkbook

# The "Spain" tab is the template for the new "Turkey" tab - duplicate it
# and place the copy right after Spain (i.e. as the new last tab).
$spain = $wb.Worksheets.Item("Spain")
$spain.Copy([System.Reflection.Missing]::Value, $spain)

$turkey = $wb.Worksheets.Item($wb.Worksheets.Count)
$turkey.Name = "Turkey"

# Market name / reference code for the new market.
# (Set the reference code first so new shared-string entries are appended
# in the same order as the authoritative edit: code, then market name.)
$turkey.Range("B4").Value = "NGC-3191/T3291"
$turkey.Range("B2").Value = "Turkey Market"

# Column D on the new sheet is a bit narrower than on the Spain sheet.
$turkey.Columns.Item(4).ColumnWidth = 19.75

# Row heights on the new tab settle back to the sheet's natural auto-fit
# (rows 3 & 5 go back to the default height, row 4 - which holds the 12pt
# reference-code font - settles at 15.6) instead of the 28.8 inherited from
# the Spain copy.
$turkey.Rows.Item(3).AutoFit()
$turkey.Rows.Item(5).AutoFit()
$turkey.Rows.Item(4).RowHeight = 15.6

# Selection / active cell bookkeeping: Turkey becomes the active tab with
# H14 selected, Spain keeps the data range selected (no longer the active tab).
$spain.Select() | Out-Null
$spain.Range("A1:D15").Select() | Out-Null

$turkey.Select() | Out-Null
$turkey.Range("H14").Select() | Out-Null
